$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 29 de Junio de 2020 a las 02:53"

# Estados Unidos (row 4) - totals refreshed
$ws.Range("B4").Value = 2637072
$ws.Range("C4").Value = 40535
$ws.Range("D4").Value = 1093456
$ws.Range("E4").Value = 1415179
$ws.Range("G4").Value = 285
$ws.Range("H4").Value = 128437

# Brasil (row 5)
$ws.Range("B5").Value = 1345254
$ws.Range("C5").Value = 29313
$ws.Range("E5").Value = 553748
$ws.Range("G5").Value = 555
$ws.Range("H5").Value = 57658

# Canada (row 22)
$ws.Range("B22").Value = 103250
$ws.Range("C22").Value = 218
$ws.Range("D22").Value = 66191
$ws.Range("E22").Value = 28537

# Argentina (row 30)
$ws.Range("B30").Value = 59933
$ws.Range("C30").Value = 2189
$ws.Range("E30").Value = 38567
$ws.Range("G30").Value = 25
$ws.Range("H30").Value = 1232

# Panama overtakes Suiza / Republica Dominicana / Afganistan / Bolivia,
# shifting those four rows down by one (rows 43-47)
$ws.Range("A43").Value = "Panama"
$ws.Range("B43").Value = 31686
$ws.Range("C43").Value = 1028
$ws.Range("D43").Value = 15470
$ws.Range("E43").Value = 15612
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 12
$ws.Range("H43").Value = 604

$ws.Range("A44").Value = "Suiza"
$ws.Range("B44").Value = 31617
$ws.Range("C44").Value = 62
$ws.Range("D44").Value = 29100
$ws.Range("E44").Value = 555
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 1962

$ws.Range("A45").Value = "Republica Dominicana"
$ws.Range("B45").Value = 31373
$ws.Range("C45").Value = 754
$ws.Range("D45").Value = 17142
$ws.Range("E45").Value = 13505
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 8
$ws.Range("H45").Value = 726

$ws.Range("A46").Value = "Afganistan"
$ws.Range("B46").Value = 30967
$ws.Range("C46").Value = 351
$ws.Range("D46").Value = 12604
$ws.Range("E46").Value = 17642
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 18
$ws.Range("H46").Value = 721

$ws.Range("A47").Value = "Bolivia"
$ws.Range("B47").Value = 30676
$ws.Range("C47").Value = 1253
$ws.Range("D47").Value = 8158
$ws.Range("E47").Value = 21548
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 36
$ws.Range("H47").Value = 970

# Venezuela overtakes Gabon, shifting Gabon down by one (rows 86-87)
$ws.Range("A86").Value = "Venezuela"
$ws.Range("B86").Value = 5297
$ws.Range("C86").Value = 167
$ws.Range("D86").Value = 1649
$ws.Range("E86").Value = 3604
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = 44

$ws.Range("A87").Value = "Gabon"
$ws.Range("B87").Value = 5209
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 2327
$ws.Range("E87").Value = 2842
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 40

# Libia (row 145)
$ws.Range("B145").Value = 762
$ws.Range("C145").Value = 35
$ws.Range("D145").Value = 196
$ws.Range("E145").Value = 545
$ws.Range("G145").Value = 3
$ws.Range("H145").Value = 21

# Polinesia Francesa (row 188)
$ws.Range("B188").Value = 62
$ws.Range("C188").Value = 2
$ws.Range("E188").Value = 2
